# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (with fresh fund-holding figures) ahead
# of the existing "2022-Q3" sheet, and rolls the quarterly summary on the
# "总计" sheet forward so it now lists 2022-Q4 / 2022-Q3 / 2021-Q4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (keeps the
#    header row / styles / column layout identical) and slot it in right
#    before "2022-Q3" so the tab order becomes:
#    总计, 2022-Q4, 2022-Q3, 2021-Q4
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Overwrite the fund rows on the new "2022-Q4" sheet with the latest
#    quarter's figures. Columns B and D:G hold numeric-looking text
#    (fund codes / figures stored as text in the source data) so force a
#    text number format before assigning, otherwise Excel would silently
#    coerce them to numbers and fund codes like "007844" would lose their
#    leading zero.
# ---------------------------------------------------------------------
function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    # Drop the temporary text number-format again now that the text value
    # is locked in, so the cell keeps the sheet's plain/default styling.
    $rng.ClearFormats()
}

# Row 2 - 001481 / 华宝油气（QDII）美元
Set-TextValue $q4 "D2" "50.10"
Set-TextValue $q4 "E2" "94.65"
Set-TextValue $q4 "F2" "2.43"
Set-TextValue $q4 "G2" "1.2174"

# Row 3 - 007844 / 华宝油气（QDII）人民币 C
Set-TextValue $q4 "B3" "007844"
$q4.Range("C3").Value = "华宝油气（QDII）人民币 C"
Set-TextValue $q4 "D3" "27.91"
Set-TextValue $q4 "E3" "94.65"
Set-TextValue $q4 "F3" "2.43"
Set-TextValue $q4 "G3" "0.6782"

# Row 4 - 162411 / 华宝油气（QDII）人民币A
Set-TextValue $q4 "B4" "162411"
$q4.Range("C4").Value = "华宝油气（QDII）人民币A"
Set-TextValue $q4 "D4" "22.19"
Set-TextValue $q4 "E4" "94.65"
Set-TextValue $q4 "F4" "2.43"
Set-TextValue $q4 "G4" "0.5392"

# ---------------------------------------------------------------------
# 3. Roll the "总计" (summary) sheet forward: 2022-Q4 becomes the newest
#    row, the old 2022-Q3 row slides to row 3, and the old 2021-Q4 row
#    slides to row 4. Row 4 is brand new, so first clone row 3's cell
#    formatting (bold/centered "A" style) onto it before filling values.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 1.63

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 2.12

$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 2.43
